# ------------------------------------------------------------------
# Applies the "icd_compare/mapping_template.xlsx" update:
#   - Columns sheet: replace the sample mapping rows (Category/Subcategory/
#     ID/Value) with a trimmed ID/Value example, add list-style data
#     validation to is_key / fill_down / confirmed_right_column.
#   - Instructions sheet: add step-by-step usage instructions.
#   - New hidden ValidationLists sheet backing the confirmed_right_column
#     drop-down.
# ------------------------------------------------------------------

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------
# 1) Columns sheet
# ---------------------------------------------------------------
$wsColumns = $wb.Worksheets.Item("Columns")

# Drop the old sample rows 4 and 5 (ID / Value) - rows 2 and 3 (Category /
# Subcategory) are overwritten in place below, becoming the new ID / Value
# sample rows.
$wsColumns.Rows("4:5").Delete()

$wsColumns.Range("A2:C2").Value = "ID"
$wsColumns.Range("A3:C3").Value = "Value"

# Clear the old is_key / fill_down markers ("Y") for these rows.
$wsColumns.Range("D2:E3").ClearContents()

# Column widths (~30/10/15 "characters").
$wsColumns.Range("A1:C1").ColumnWidth = 30
$wsColumns.Range("D1").ColumnWidth = 10
$wsColumns.Range("E1").ColumnWidth = 15

# ---------------------------------------------------------------
# 2) ValidationLists sheet (new, hidden, added after Instructions)
# ---------------------------------------------------------------
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$wsLists = $wb.Worksheets.Add([Type]::Missing, $lastSheet)
$wsLists.Name = "ValidationLists"
$wsLists.Range("A1").Value = "ID"
$wsLists.Range("A2").Value = "Value"
$wsLists.Visible = $false

# ---------------------------------------------------------------
# 3) Data validations on the Columns sheet
#    (xlValidateList / xlValidAlertStop / xlBetween = 3 / 1 / 1)
# ---------------------------------------------------------------
$xlValidateList = [Microsoft.Office.Interop.Excel.XlDVType]::xlValidateList
$xlValidAlertStop = [Microsoft.Office.Interop.Excel.XlDVAlertStyle]::xlValidAlertStop
$xlBetween = [Microsoft.Office.Interop.Excel.XlFormatConditionOperator]::xlBetween

$wsColumns.Range("D2:D3").Validation.Add($xlValidateList, $xlValidAlertStop, $xlBetween, '"Y,N"')
$wsColumns.Range("E2:E3").Validation.Add($xlValidateList, $xlValidAlertStop, $xlBetween, '"Y,N"')
$wsColumns.Range("C2:C3").Validation.Add($xlValidateList, $xlValidAlertStop, $xlBetween, "=ValidationLists!`$A`$1:`$A`$2")

# ---------------------------------------------------------------
# 4) Instructions sheet - append the step-by-step guide
# ---------------------------------------------------------------
$wsInstructions = $wb.Worksheets.Item("Instructions")

# Write the "Step N" labels first, then the action text, so that the
# shared-string table is populated in the same order as the source edit.
$wsInstructions.Range("A2").Value = "Step 1"
$wsInstructions.Range("A3").Value = "Step 2"
$wsInstructions.Range("A4").Value = "Step 3"
$wsInstructions.Range("A5").Value = "Step 4"
$wsInstructions.Range("A6").Value = "Step 5"
$wsInstructions.Range("A7").Value = "Step 6"

$wsInstructions.Range("B2").Value = "Review the 'Columns' sheet."
$wsInstructions.Range("B3").Value = "Verify 'confirmed_right_column' matches the correct column in the Right file. Clear it if you don't want to compare that column."
$wsInstructions.Range("B4").Value = "Mark Key columns by entering 'Y' in 'is_key'. Keys are used to join rows."
$wsInstructions.Range("B5").Value = "Mark columns that need Forward Fill (e.g. parent IDs) by entering 'Y' in 'fill_down'. This fills empty cells with the value from the row above."
$wsInstructions.Range("B6").Value = "Save this workbook."
$wsInstructions.Range("B7").Value = "Run the diff script again with --mapping-confirmed."

# Re-select the Columns sheet (tabSelected="1" in the original file).
$wsColumns.Activate()
